$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, pushing existing rows 24-32 down to 25-33.
$ws.Rows.Item(24).Insert()

# Fill the new row 24 with the data for the new weekly price entry.
$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value = "Maule"
$ws.Cells.Item(24, 4).Value = 44529
$ws.Cells.Item(24, 4).Style = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat
$ws.Cells.Item(24, 5).Value = 7
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100101
$ws.Cells.Item(24, 8).Value = "Berries"
$ws.Cells.Item(24, 9).Value = 100101001
$ws.Cells.Item(24, 10).Value = "Arándano (blue)"
$ws.Cells.Item(24, 11).Value = "Sin especificar"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 3800
$ws.Cells.Item(24, 15).Value = 3800
$ws.Cells.Item(24, 16).Value = 3800
$ws.Cells.Item(24, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(24, 18).Value = "Provincia de Linares"
$ws.Cells.Item(24, 19).Value = 1900
$ws.Cells.Item(24, 20).Value = 2

$wb.Save()
